# Add two new columns to the sheet: I0 (col I) and IF (col J), mirroring
# the existing header/style conventions used by column H ("IP").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting already used on H1 (bold, centered, bordered)
# by copying H1's format onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I0 and IF columns, indexed by row number (2..28)
$i0 = @{
    2 = 1; 3 = 1; 4 = 1; 5 = 1; 6 = 1; 7 = 8; 8 = 7; 9 = 4; 10 = 3;
    11 = 1; 12 = 1; 13 = 1; 14 = 1; 15 = 1; 16 = 1; 17 = 1; 18 = 1; 19 = 1;
    20 = 1; 21 = 1; 22 = 1; 23 = 1; 24 = 1; 25 = 1; 26 = 1; 27 = 1; 28 = 1
}
$if_ = @{
    2 = 2; 3 = 5; 4 = 5; 5 = 7; 6 = 5; 7 = 8; 8 = 7; 9 = 5; 10 = 4;
    11 = 2; 12 = 5; 13 = 3; 14 = 6; 15 = 6; 16 = 3; 17 = 5; 18 = 4; 19 = 6;
    20 = 7; 21 = 6; 22 = 7; 23 = 6; 24 = 6; 25 = 6; 26 = 5; 27 = 5; 28 = 2
}

foreach ($r in 2..28) {
    $ws.Cells.Item($r, 9).Value = $i0[$r]
    $ws.Cells.Item($r, 10).Value = $if_[$r]
}
